$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Feria Lagunitas de
# Puerto Montt - Ajo". It belongs chronologically right before the data
# that is currently on row 89, so insert a fresh row there (pushing every
# row from 89 down through 206 down by one, which is exactly what the
# target dimension A1:R207 reflects) and fill it in.
$ws.Rows.Item(89).Insert()

$ws.Range("A89").Value = 4
$ws.Range("B89").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C89").Value = "Los Lagos"
$ws.Range("D89").Value = 44579
$ws.Range("E89").Value = 10
$ws.Range("F89").Value = 100112003
$ws.Range("G89").Value = "Ajo"
$ws.Range("H89").Value = "Chino"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 240
$ws.Range("K89").Value = 21000
$ws.Range("L89").Value = 22000
$ws.Range("M89").Value = 21500
$ws.Range("N89").Value = '$/caja 10 kilos'
$ws.Range("O89").Value = "China"
$ws.Range("P89").Value = 2150
$ws.Range("Q89").Value = 10
$ws.Range("R89").Value = "Hortaliza"
